# Adds a new "MergeSort" worksheet (with timing-metric data and a line chart)
# after the existing "QuickSort" sheet, mirroring the structure already used
# for "Cambio Residuo" and "QuickSort".

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet after the last existing sheet -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "MergeSort"

# --- 2. Fill in the metric values (MergeSort!A2:A51) ------------------------
$values = @(
    [double]"7.9999999999999996E-7",
    [double]"1.3E-6",
    [double]"2.0999999999999998E-6",
    [double]"2.7E-6",
    [double]"3.4999999999999999E-6",
    [double]"4.0999999999999997E-6",
    [double]"6.3999999999999997E-6",
    [double]"5.2000000000000002E-6",
    [double]"5.9000000000000003E-6",
    [double]"6.7000000000000002E-6",
    [double]"7.4000000000000003E-6",
    [double]"7.7999999999999999E-6",
    [double]"8.6999999999999997E-6",
    [double]"9.3000000000000007E-6",
    [double]"1.01E-5",
    [double]"1.0699999999999999E-5",
    [double]"1.17E-5",
    [double]"1.24E-5",
    [double]"1.27E-5",
    [double]"1.3900000000000001E-5",
    [double]"1.4600000000000001E-5",
    [double]"1.5699999999999999E-5",
    [double]"1.6399999999999999E-5",
    [double]"1.6699999999999999E-5",
    [double]"1.7900000000000001E-5",
    [double]"1.8199999999999999E-5",
    [double]"1.9199999999999999E-5",
    [double]"1.98E-5",
    [double]"2.05E-5",
    [double]"2.09E-5",
    [double]"2.6599999999999999E-5",
    [double]"2.8600000000000001E-5",
    [double]"3.1199999999999999E-5",
    [double]"3.26E-5",
    [double]"2.6100000000000001E-5",
    [double]"2.8900000000000001E-5",
    [double]"2.9799999999999999E-5",
    [double]"3.0800000000000003E-5",
    [double]"2.9E-5",
    [double]"3.2499999999999997E-5",
    [double]"3.1099999999999997E-5",
    [double]"3.7499999999999997E-5",
    [double]"3.2499999999999997E-5",
    [double]"3.3000000000000003E-5",
    [double]"3.3800000000000002E-5",
    [double]"3.4499999999999998E-5",
    [double]"3.4900000000000001E-5",
    [double]"3.6000000000000001E-5",
    [double]"3.6300000000000001E-5",
    [double]"3.79E-5"
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}

$lastRow = $startRow + $values.Count - 1
$dataRange = $ws.Range("A$($startRow):A$($lastRow)")

# Scientific notation number format (matches the new style used in the file)
$dataRange.NumberFormat = "0.00E+00"

# --- 3. Add a line chart plotting the metric values -------------------------
$co = $ws.ChartObjects().Add(107.19, 19.5, 350.625, 216)
$chart = $co.Chart
$chart.ChartType = 4
$chart.SetSourceData($dataRange)
$chart.HasLegend = $false
$chart.HasTitle = $true

# --- 4. Make the new sheet the active tab (as in the edited workbook) ------
$ws.Activate()
$ws.Range("G22").Select()
